# Fixed some bugs in genStackReels
# Update reel-stack table values (rows 2-21) to corrected data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 101
$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 15
$ws.Range("F2").Value = 15
$ws.Range("A4").Value = 401
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 48
$ws.Range("D4").Value = 67
$ws.Range("E4").Value = 75
$ws.Range("F4").Value = 45
$ws.Range("A5").Value = 901
$ws.Range("B5").Value = 16
$ws.Range("C5").Value = 15
$ws.Range("D5").Value = 45
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 60
$ws.Range("A6").Value = 1201
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 10
$ws.Range("A7").Value = 1203
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = 15
$ws.Range("D7").Value = 15
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 15
$ws.Range("A8").Value = 902
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("A9").Value = 501
$ws.Range("B9").Value = 9
$ws.Range("C9").Value = 52
$ws.Range("D9").Value = 30
$ws.Range("E9").Value = 75
$ws.Range("F9").Value = 45
$ws.Range("A10").Value = 701
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 90
$ws.Range("D10").Value = 45
$ws.Range("E10").Value = 97
$ws.Range("F10").Value = 15
$ws.Range("A11").Value = 601
$ws.Range("C11").Value = 60
$ws.Range("D11").Value = 67
$ws.Range("F11").Value = 42
$ws.Range("A12").Value = 801
$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 67
$ws.Range("D12").Value = 65
$ws.Range("E12").Value = 52
$ws.Range("F12").Value = 45
$ws.Range("A13").Value = 1001
$ws.Range("B13").Value = 18
$ws.Range("C13").Value = 30
$ws.Range("D13").Value = 75
$ws.Range("E13").Value = 60
$ws.Range("F13").Value = 72
$ws.Range("A15").Value = 1202
$ws.Range("B15").Value = 2
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 10
$ws.Range("A16").Value = 1101
$ws.Range("C16").Value = 15
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 30
$ws.Range("F16").Value = 0
$ws.Range("A18").Value = 3
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 3
$ws.Range("F18").Value = 3
$ws.Range("A19").Value = 1
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 2
$ws.Range("A20").Value = 502
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = 0
$ws.Range("A21").Value = 802
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 5
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = 0
